{"js": "// Add four new TODO bullets to the end of the \"To do 26/08/14\" numbered\n// list (the list that uses numId=2 / style \"ListParagraph\"). The last\n// paragraph of the document is an (empty) list item that carries the\n// \"_GoBack\" bookmark; we fill it with the final new bullet's text and\n// insert three brand-new list paragraphs with the other bullets right\n// before it, so the bookmark stays attached to the very last paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the document - currently empty, holds the\n// \"_GoBack\" bookmark.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert the first three new bullets before that paragraph, preserving\n// its (ListParagraph / numId=2) formatting.\nconst bullet1 = lastParagraph.insertParagraph(\n  \"[calculation] recalculate the difference between expected (rouse chain) model encounter probability and the observed to spot peaks of the experimental data\",\n  \"Before\"\n);\nawait context.sync();\n\nconst bullet2 = bullet1.insertParagraph(\n  \"Present the two-sided bead encounter frequency and show that by dividing each bead\\u2019s encounter frequency by the sum of encounters the symmetry is broken \",\n  \"After\"\n);\nawait context.sync();\n\nbullet2.insertParagraph(\n  \"Finish listing the peaks of the encounter frequencies in the experimental data\",\n  \"After\"\n);\nawait context.sync();\n\n// Put the final bullet's text into the original last (bookmark) paragraph\n// so the \"_GoBack\" bookmark remains on the document's last paragraph.\nlastParagraph.insertText(\n  \"Calculate the beta values of the experimental data with no peaks (remove peaks by assigning the peaks with the neighbors encounter values)\",\n  \"Start\"\n);\nawait context.sync();\n", "ps1": "# Add four new TODO bullets to the end of the \"To do 26/08/14\" numbered\n# list (numId=2 / style \"ListParagraph\"). The last paragraph of the\n# document is an (empty) list item that carries the \"_GoBack\" bookmark;\n# we insert three brand-new list paragraphs (inheriting the same list\n# formatting) right before it, fill in the first three bullets' text,\n# and put the fourth bullet's text into the original last paragraph so\n# the bookmark stays attached to the document's very last paragraph.\n\n$d = $word.ActiveDocument\n$n = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($n)\n$lastRange = $lastParagraph.Range\n\n$lastRange.InsertParagraphBefore()\n$lastRange.InsertParagraphBefore()\n$lastRange.InsertParagraphBefore()\n\n$n2 = $d.Paragraphs.Count\n$bullet1 = $d.Paragraphs.Item($n2 - 3)\n$bullet2 = $d.Paragraphs.Item($n2 - 2)\n$bullet3 = $d.Paragraphs.Item($n2 - 1)\n$bullet4 = $d.Paragraphs.Item($n2)\n\n$bullet1.Range.Text = \"[calculation] recalculate the difference between expected (rouse chain) model encounter probability and the observed to spot peaks of the experimental data\"\n$bullet2.Range.Text = \"Present the two-sided bead encounter frequency and show that by dividing each bead\u2019s encounter frequency by the sum of encounters the symmetry is broken \"\n$bullet3.Range.Text = \"Finish listing the peaks of the encounter frequencies in the experimental data\"\n$bullet4.Range.Text = \"Calculate the beta values of the experimental data with no peaks (remove peaks by assigning the peaks with the neighbors encounter values)\"\n"}
